$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25. This shifts the old rows 25-29
# (years 2017-2021) down to rows 26-30, opening up a blank row 25
# for the new year-2016 estimates.
$ws.Rows("25:25").Insert()

# The newly-inserted row's A cell does not pick up the same formatting
# as the rest of column A, so copy it over from the row above (A24),
# which carries the correct style.
$ws.Range("A24").Copy($ws.Range("A25"))

# Column A is a simple running index (row number - 2) that is independent
# of the data rows, so restore/continue that sequence for row 25 (new)
# through row 30 (the old row 29, now shifted down one).
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25
$ws.Range("A28").Value = 26
$ws.Range("A29").Value = 27
$ws.Range("A30").Value = 28

# Fill in the new row 25 with the year-2016 estimates.
$ws.Range("B25").Value = 2016
$ws.Range("C25").Value = 315821
$ws.Range("D25").Value = 70790.49112339407
$ws.Range("E25").Value = 1.516615146580414
$ws.Range("F25").Value = 118245
$ws.Range("G25").Value = 2257636
$ws.Range("H25").Value = 0.9159999999999999
$ws.Range("I25").Value = "Total Manufactura sin ENGE"
$ws.Range("J25").Value = 47454.50887660593
$ws.Range("K25").Value = 1941815
$ws.Range("L25").Value = 0.779838252556603
$ws.Range("M25").Value = 1.305629871784614
$ws.Range("N25").Value = 0.5972889173335725
$ws.Range("O25").Value = 0.05237558224620798
$ws.Range("P25").Value = 0.2241475111642166
$ws.Range("Q25").Value = 0.02443822345414261
$ws.Range("R25").Value = 4.279618508306799
$ws.Range("S25").Value = 0.4665957380533359
$ws.Range("T25").Value = 0.5986764017370212
$ws.Range("U25").Value = 0.1398901328646425
